# Data provider is added

$wb = $excel.ActiveWorkbook

# Update LoginTest credentials
$loginSheet = $wb.Worksheets.Item("LoginTest")
$loginSheet.Range("A2").Value = "leoalak@gmail.com"
$loginSheet.Range("B2").Value = "Toma*1996"
[void]$loginSheet.Range("O21").Select()

# Add a new worksheet "searchSomething" positioned between LoginTest and CreateAccountTest
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "searchSomething"
$newSheet.Range("A1").Value = "key"
$newSheet.Range("A2").Value = "alak"
[void]$newSheet.Range("B8").Select()

# Update CreateAccountTest data
$createAccountSheet = $wb.Worksheets.Item("CreateAccountTest")
$createAccountSheet.Range("A2").Value = "Alak"
[void]$createAccountSheet.Range("M7").Select()

# Restore focus on the newly inserted sheet (matches activeTab position)
[void]$newSheet.Select()
